$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. New mirrored table in columns G:K - "FRONT IR SENSOR; After Weird
#    Behavior" - tracking the same trial numbers as the front-sensor
#    table in A:E, but with the trial/avg cells still blank.
# ------------------------------------------------------------------

# Header (merged G1:K1), copy format from the existing A1:E1 banner.
$ws.Range("A1:E1").Copy($ws.Range("G1")) | Out-Null
$ws.Range("G1").Value = "FRONT IR SENSOR; After Weird Behavior"

# Column headers row (Distance / Trial 1 / Trial 2 / Trial 3 / AVG).
$ws.Range("A2:E2").Copy($ws.Range("G2")) | Out-Null

# Distance values + blank (but formatted) trial/avg cells for rows 3-24.
$ws.Range("A3:E24").Copy($ws.Range("G3")) | Out-Null
$ws.Range("H3:K24").ClearContents() | Out-Null

# ------------------------------------------------------------------
# 2. Move the chart to the right of the new table.
# ------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$co.Left = 1000.3065618848425
$co.Top = 9

# ------------------------------------------------------------------
# 3. Scroll back to the top and move the active selection.
# ------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A1"))
$ws.Range("N23").Select() | Out-Null

Write-Host "done"
